$p = $ppt.ActivePresentation

$newDate = "11/16/2025"
$ppPlaceholderDate = 16

function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDatePh = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDatePh = $true
            }
        } catch {
            $isDatePh = $false
        }
        if ($isDatePh -and $sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -ne $newDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Update the slide master's date placeholder.
$master = $p.SlideMaster
Update-DateShape $master.Shapes

# Update every slide layout's date placeholder.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateShape $layout.Shapes
}
